$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O39").Value = 40
$ws.Range("P39").Value = 628
$ws.Range("Q39").Value = 388
$ws.Range("R39").Value = 184
$ws.Range("S39").Value = 0.2
$ws.Range("T39").Value = 0.068
$ws.Range("U39").Value = 0.22
$ws.Range("V39").Value = 0.223
$ws.Range("W39").Value = 0.52
$ws.Range("O45").Value = 27
$ws.Range("P45").Value = 91
$ws.Range("Q45").Value = 23
$ws.Range("R45").Value = 51
$ws.Range("S45").Value = 4.8
$ws.Range("U45").Value = 0.058
$ws.Range("V45").Value = 0.181
$ws.Range("W45").Value = 0.516
$ws.Range("X45").Value = 0.256
$ws.Range("O80").Value = 29
$ws.Range("P80").Value = 134
$ws.Range("Q80").Value = 128
$ws.Range("S80").Value = -11.8
$ws.Range("T80").Value = 0.133
$ws.Range("U80").Value = 0.242
$ws.Range("V80").Value = 0.168
$ws.Range("W80").Value = 0.599
$ws.Range("X80").Value = 0.06
$ws.Range("O85").Value = 41
$ws.Range("P85").Value = 393
$ws.Range("Q85").Value = 173
$ws.Range("R85").Value = 338
$ws.Range("S85").Value = 3.4
$ws.Range("T85").Value = 0.016
$ws.Range("V85").Value = 0.141
$ws.Range("W85").Value = 0.58
$ws.Range("X85").Value = 0.378
$ws.Range("O139").Value = 26
$ws.Range("P139").Value = 423
$ws.Range("Q139").Value = 98
$ws.Range("R139").Value = 81
$ws.Range("S139").Value = -11
$ws.Range("T139").Value = 0.024
$ws.Range("V139").Value = 0.233
$ws.Range("X139").Value = 0.17
$ws.Range("O153").Value = 39
$ws.Range("P153").Value = 418
$ws.Range("Q153").Value = 91
$ws.Range("R153").Value = 105
$ws.Range("S153").Value = 1.9
$ws.Range("T153").Value = 0.014
$ws.Range("U153").Value = 0.08799999999999999
$ws.Range("W153").Value = 0.572
$ws.Range("X153").Value = 0.172
$ws.Range("O179").Value = 41
$ws.Range("P179").Value = 457
$ws.Range("Q179").Value = 156
$ws.Range("R179").Value = 70
$ws.Range("S179").Value = 0.1
$ws.Range("U179").Value = 0.079
$ws.Range("W179").Value = 0.623
$ws.Range("X179").Value = 0.08500000000000001
$ws.Range("O181").Value = 39
$ws.Range("P181").Value = 273
$ws.Range("Q181").Value = 129
$ws.Range("S181").Value = -2.4
$ws.Range("T181").Value = 0.035
$ws.Range("V181").Value = 0.112
$ws.Range("W181").Value = 0.611
$ws.Range("X181").Value = 0.07000000000000001
$ws.Range("O209").Value = 37
$ws.Range("P209").Value = 356
$ws.Range("Q209").Value = 176
$ws.Range("R209").Value = 100
$ws.Range("S209").Value = -6
$ws.Range("T209").Value = 0.052
$ws.Range("U209").Value = 0.141
$ws.Range("V209").Value = 0.194
$ws.Range("W209").Value = 0.533
$ws.Range("X209").Value = 0.176
$ws.Range("O256").Value = 24
$ws.Range("P256").Value = 413
$ws.Range("Q256").Value = 127
$ws.Range("R256").Value = 112
$ws.Range("S256").Value = 4.4
$ws.Range("T256").Value = 0.07199999999999999
$ws.Range("U256").Value = 0.102
$ws.Range("V256").Value = 0.202
$ws.Range("W256").Value = 0.639
$ws.Range("X256").Value = 0.229
$ws.Range("O276").Value = 6
$ws.Range("S276").Value = -5.4
$ws.Range("U276").Value = 0.083
$ws.Range("V276").Value = 0.192
$ws.Range("W276").Value = 0.425
$ws.Range("X276").Value = 0.211
$ws.Range("O284").Value = 41
$ws.Range("P284").Value = 464
$ws.Range("Q284").Value = 180
$ws.Range("R284").Value = 51
$ws.Range("S284").Value = -0.3
$ws.Range("T284").Value = 0.028
$ws.Range("U284").Value = 0.123
$ws.Range("V284").Value = 0.167
$ws.Range("W284").Value = 0.576
$ws.Range("X284").Value = 0.065
$ws.Range("O300").Value = 26
$ws.Range("P300").Value = 189
$ws.Range("Q300").Value = 104
$ws.Range("R300").Value = 15
$ws.Range("S300").Value = 5.5
$ws.Range("T300").Value = 0.076
$ws.Range("U300").Value = 0.185
$ws.Range("V300").Value = 0.189
$ws.Range("W300").Value = 0.643
$ws.Range("X300").Value = 0.064
$ws.Range("O301").Value = 38
$ws.Range("P301").Value = 468
$ws.Range("Q301").Value = 189
$ws.Range("S301").Value = -7.7
$ws.Range("V301").Value = 0.203
$ws.Range("W301").Value = 0.551
$ws.Range("X301").Value = 0.095
$ws.Range("O308").Value = 9
$ws.Range("Q308").Value = 15
$ws.Range("R308").Value = 3
$ws.Range("S308").Value = 2.6
$ws.Range("T308").Value = 0.027
$ws.Range("U308").Value = 0.212
$ws.Range("X308").Value = 0.079
$ws.Range("O346").Value = 26
$ws.Range("S346").Value = -14.1
$ws.Range("T346").Value = 0.01
$ws.Range("U346").Value = 0.093
$ws.Range("V346").Value = 0.221
$ws.Range("W346").Value = 0.552
$ws.Range("X346").Value = 0.156
$ws.Range("O389").Value = 32
$ws.Range("P389").Value = 325
$ws.Range("Q389").Value = 124
$ws.Range("R389").Value = 83
$ws.Range("S389").Value = 2.6
$ws.Range("U389").Value = 0.137
$ws.Range("V389").Value = 0.183
$ws.Range("W389").Value = 0.591
$ws.Range("X389").Value = 0.173
$ws.Range("O416").Value = 25
$ws.Range("P416").Value = 102
$ws.Range("S416").Value = 5.5
$ws.Range("T416").Value = 0.033
$ws.Range("V416").Value = 0.135
$ws.Range("W416").Value = 0.5610000000000001
$ws.Range("O444").Value = 25
$ws.Range("P444").Value = 111
$ws.Range("Q444").Value = 62
$ws.Range("R444").Value = 10
$ws.Range("S444").Value = 1.2
$ws.Range("T444").Value = 0.06900000000000001
$ws.Range("U444").Value = 0.203
$ws.Range("W444").Value = 0.705
$ws.Range("X444").Value = 0.07199999999999999
$ws.Range("O453").Value = 14
$ws.Range("P453").Value = 19
$ws.Range("Q453").Value = 9
$ws.Range("S453").Value = -2.1
$ws.Range("T453").Value = 0.043
$ws.Range("U453").Value = 0.127
$ws.Range("V453").Value = 0.171
$ws.Range("W453").Value = 0.5669999999999999
$ws.Range("X453").Value = 0.222
$ws.Range("O459").Value = 40
$ws.Range("P459").Value = 465
$ws.Range("Q459").Value = 101
$ws.Range("R459").Value = 143
$ws.Range("S459").Value = -1.7
$ws.Range("T459").Value = 0.028
$ws.Range("U459").Value = 0.066
$ws.Range("V459").Value = 0.226
$ws.Range("W459").Value = 0.494
$ws.Range("X459").Value = 0.215
$ws.Range("O471").Value = 37
$ws.Range("P471").Value = 446
$ws.Range("Q471").Value = 152
$ws.Range("R471").Value = 108
$ws.Range("S471").Value = -1.2
$ws.Range("U471").Value = 0.113
$ws.Range("V471").Value = 0.197
$ws.Range("W471").Value = 0.518
$ws.Range("X471").Value = 0.154
$ws.Range("O481").Value = 23
$ws.Range("P481").Value = 85
$ws.Range("Q481").Value = 50
$ws.Range("R481").Value = 85
$ws.Range("S481").Value = 3
$ws.Range("T481").Value = 0.03
$ws.Range("U481").Value = 0.089
$ws.Range("V481").Value = 0.105
$ws.Range("W481").Value = 0.533
$ws.Range("X481").Value = 0.304
$ws.Range("O496").Value = 41
$ws.Range("P496").Value = 986
$ws.Range("Q496").Value = 233
$ws.Range("R496").Value = 207
$ws.Range("S496").Value = 0
$ws.Range("U496").Value = 0.15
$ws.Range("V496").Value = 0.26
$ws.Range("W496").Value = 0.626
$ws.Range("O505").Value = 36
$ws.Range("P505").Value = 880
$ws.Range("Q505").Value = 387
$ws.Range("R505").Value = 133
$ws.Range("S505").Value = 3.7
$ws.Range("V505").Value = 0.309
$ws.Range("W505").Value = 0.592
$ws.Range("O516").Value = 32
$ws.Range("P516").Value = 160
$ws.Range("Q516").Value = 97
$ws.Range("R516").Value = 49
$ws.Range("S516").Value = -7.5
$ws.Range("T516").Value = 0.083
$ws.Range("V516").Value = 0.169
$ws.Range("W516").Value = 0.5580000000000001
$ws.Range("X516").Value = 0.174
